$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.823.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.996.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.51%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "359.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.637"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.36%  "

$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.63%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.462.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.073.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.22%  "

$ws.Range("E17").Value = "  +1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.799.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "272.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.182"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.08%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +7.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.07%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.94%  "

$ws.Range("E42").Value = "  +3.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.159.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.52%  "

$ws.Range("E48").Value = "  -4.91%  "

$ws.Range("B49").Value = "BEAM"
$ws.Range("C49").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0356"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.59%  "

$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.945"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
